# Add a new test user record as row 17 on the "info" sheet (sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")

$ws.Range("A17").Value = "Emma"
$ws.Range("B17").Value = "Watson"
$ws.Range("C17").Value = "Test_User1"
$ws.Range("D17").Value = "Test_User1"
$ws.Range("E17").Value = "93/4 Khawsan Rd."
$ws.Range("F17").Value = "099-999-9999"
